$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("D2").Value = 0.0303
$ws.Range("E2").Value = 0.1458
$ws.Range("F2").Value = 219.96
$ws.Range("G2").Value = 0.2078189300411522
$ws.Range("H2").Value = 0.0183
$ws.Range("I2").Value = -1.112

# Remove row 3 entirely (PA state row)
$ws.Range("A3:I3").EntireRow.Delete()
